$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("labels")

# New rows appended to the "labels" sheet (rows 22-30), matching the
# "Case study from LEAP trial" / "Chen, et al paper on maternal
# investment" / famuss additions from the Jul 15 commit.
$rows = @(
    @("basicExampleOfPeanutAllergies", "section", "Case study from LEAP trial"),
    @("peanutStudyResultsDF", "table", "results for 5 studies from LEAP"),
    @("peanutStudyResults", "table", "summary of LEAP outcome"),
    @("basicExampleOfFrogAltitude", "section", "Chen, et al paper on maternal investment"),
    @("FrogAltitudeDF", "table", "sample data matrix for frog altitude data"),
    @("famuss_height_weight", "figure", "scatterplot of height vs weight"),
    @("famuss_height_bmi", "figure", "scatterplot of height vs bmi"),
    @("FAMuSSDF", "table", "table of 4 cases from famuss"),
    @("FAMuSS_subset_Variables", "table", "definitions in famuss")
)

$r = 22
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

$ws.Range("C30").Select()
